{"js": "// The paragraph containing the \"<id>\" tag currently has the tag text\n// split across three runs: \"<id>\" (Courier New / #7f6000), \"p069r_2\"\n// (default formatting) and \"</id>\" (Courier New / #7f6000). Merge them\n// into a single run containing the full \"<id>p069r_2</id>\" text, using\n// the formatting of the first run.\nconst results = context.document.body.search(\"<id>p069r_2</id>\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the split <id>p069r_2</id> run sequence.\");\n}\n\nconst target = results.items[0];\n// Replacing the whole matched range with its own text collapses the\n// three runs into one run (taking on the formatting of the first run\n// in the range), exactly mirroring the OOXML merge in the diff.\ntarget.insertText(\"<id>p069r_2</id>\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The paragraph containing the \"<id>\" tag currently has its text split\n# across three runs: \"<id>\" (Courier New / #7f6000), \"p069r_2\" (default\n# formatting) and \"</id>\" (Courier New / #7f6000). Merge them into a\n# single run containing \"<id>p069r_2</id>\", keeping the formatting of\n# the first run (\"<id>\").\n\n$searchRng = $d.Content\n$find = $searchRng.Find\n$find.ClearFormatting()\n$find.Text = \"<id>\"\n$found = $find.Execute()\n\nif ($found) {\n    # $searchRng has been collapsed by Execute() to the matched \"<id>\" text.\n    $idEnd = $searchRng.End\n    $nextChar = $d.Range($idEnd, $idEnd + 1)\n\n    # Only merge if the text right after \"<id>\" is not already in a run\n    # with the same formatting (i.e. skip if already merged).\n    $alreadyMerged = ($nextChar.Font.Name -eq $searchRng.Font.Name) -and `\n        ($nextChar.Font.Color -eq $searchRng.Font.Color) -and `\n        ($nextChar.Font.Size -eq $searchRng.Font.Size)\n\n    if (-not $alreadyMerged) {\n        $tailText = \"p069r_2</id>\"\n        $tail = $d.Range($idEnd, $idEnd + $tailText.Length)\n\n        if ($tail.Text -eq $tailText) {\n            # Remove the trailing runs' text, then append the same text\n            # right after \"<id>\" so it becomes part of the first run (and\n            # picks up its formatting) -- matching the merge in the diff.\n            $tail.Delete()\n            $searchRng.InsertAfter($tailText)\n        }\n    }\n}\n"}
